# Add a new "Save" column (H) to the s_vals sheet, mirroring the
# header style used by the existing "sum" column (G), and fill the
# data rows with placeholder zero values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font, thin border, centered/top
# alignment) from G1 onto the new header cell H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New numeric data cells for the "Save" column.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
